$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 <= source row 12 (Id 112044190)
$ws.Range("A10").Value = 112044190
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("AF10").ClearContents()
$ws.Range("Q10").Value = 554682
$ws.Range("R10").Value = 6698694
$ws.Range("Z10").ClearContents()
$ws.Range("AB10").ClearContents()

# Row 11 <= source row 10 (Id 112044194)
$ws.Range("A11").Value = 112044194
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("AF11").ClearContents()
$ws.Range("Q11").Value = 554746
$ws.Range("R11").Value = 6698619
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()

# Row 12 <= source row 18 (Id 112044187)
$ws.Range("A12").Value = 112044187
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("J12").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("AF12").ClearContents()
$ws.Range("Q12").Value = 554629
$ws.Range("R12").Value = 6698775
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()

# Row 13 <= source row 21 (Id 112044170)
$ws.Range("A13").Value = 112044170
$ws.Range("B13").Value = 89845
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 1209
$ws.Range("F13").Value = "Rynkskinn"
$ws.Range("G13").Value = "Phlebia centrifuga"
$ws.Range("H13").Value = "P.Karst."
$ws.Range("J13").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("AF13").ClearContents()
$ws.Range("Q13").Value = 554745
$ws.Range("R13").Value = 6698641
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()

# Row 14 <= source row 15 (Id 112044174)
$ws.Range("A14").Value = 112044174
$ws.Range("B14").Value = 96348
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("J14").ClearContents()
$ws.Range("K14").Value = "överblommad"
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("AF14").ClearContents()
$ws.Range("Q14").Value = 554690
$ws.Range("R14").Value = 6698722
$ws.Range("Z14").ClearContents()
$ws.Range("AB14").ClearContents()

# Row 15 <= source row 14 (Id 112044188)
$ws.Range("A15").Value = 112044188
$ws.Range("B15").Value = 96348
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = "Knärot"
$ws.Range("G15").Value = "Goodyera repens"
$ws.Range("H15").Value = "(L.) R. Br."
$ws.Range("J15").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("AF15").ClearContents()
$ws.Range("Q15").Value = 554647
$ws.Range("R15").Value = 6698760
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()

# Row 16 <= source row 20 (Id 112044163)
$ws.Range("A16").Value = 112044163
$ws.Range("B16").Value = 56543
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 103021
$ws.Range("F16").Value = "Talltita"
$ws.Range("G16").Value = "Poecile montanus"
$ws.Range("H16").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("J16").ClearContents()
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = "spel/sång"
$ws.Range("N16").ClearContents()
$ws.Range("AF16").ClearContents()
$ws.Range("Q16").Value = 554650
$ws.Range("R16").Value = 6698762
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()

# Row 17 <= source row 16 (Id 112044186)
$ws.Range("A17").Value = 112044186
$ws.Range("B17").Value = 96348
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("AF17").ClearContents()
$ws.Range("Q17").Value = 554675
$ws.Range("R17").Value = 6698785
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()

# Row 18 <= source row 13 (Id 112044191)
$ws.Range("A18").Value = 112044191
$ws.Range("B18").Value = 96348
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
$ws.Range("J18").ClearContents()
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("AF18").ClearContents()
$ws.Range("Q18").Value = 554719
$ws.Range("R18").Value = 6698669
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()

# Row 19 <= source row 11 (Id 112044189)
$ws.Range("A19").Value = 112044189
$ws.Range("B19").Value = 96348
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = "Knärot"
$ws.Range("G19").Value = "Goodyera repens"
$ws.Range("H19").Value = "(L.) R. Br."
$ws.Range("J19").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("AF19").ClearContents()
$ws.Range("Q19").Value = 554686
$ws.Range("R19").Value = 6698721
$ws.Range("Z19").ClearContents()
$ws.Range("AB19").ClearContents()

# Row 20 <= source row 17 (Id 112044185)
$ws.Range("A20").Value = 112044185
$ws.Range("B20").Value = 96348
$ws.Range("D20").Value = "VU"
$ws.Range("E20").Value = 220787
$ws.Range("F20").Value = "Knärot"
$ws.Range("G20").Value = "Goodyera repens"
$ws.Range("H20").Value = "(L.) R. Br."
$ws.Range("J20").ClearContents()
$ws.Range("K20").ClearContents()
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("AF20").ClearContents()
$ws.Range("Q20").Value = 554752
$ws.Range("R20").Value = 6698637
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()

# Row 21 <= source row 22 (Id 112044192)
$ws.Range("A21").Value = 112044192
$ws.Range("B21").Value = 96348
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 220787
$ws.Range("F21").Value = "Knärot"
$ws.Range("G21").Value = "Goodyera repens"
$ws.Range("H21").Value = "(L.) R. Br."
$ws.Range("J21").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("AF21").ClearContents()
$ws.Range("Q21").Value = 554727
$ws.Range("R21").Value = 6698622
$ws.Range("Z21").ClearContents()
$ws.Range("AB21").ClearContents()

# Row 22 <= source row 19 (Id 112044193)
$ws.Range("A22").Value = 112044193
$ws.Range("B22").Value = 96348
$ws.Range("D22").Value = "VU"
$ws.Range("E22").Value = 220787
$ws.Range("F22").Value = "Knärot"
$ws.Range("G22").Value = "Goodyera repens"
$ws.Range("H22").Value = "(L.) R. Br."
$ws.Range("J22").ClearContents()
$ws.Range("K22").ClearContents()
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("AF22").ClearContents()
$ws.Range("Q22").Value = 554737
$ws.Range("R22").Value = 6698616
$ws.Range("Z22").ClearContents()
$ws.Range("AB22").ClearContents()

# Row 24 <= source row 28 (Id 112044162)
$ws.Range("A24").Value = 112044162
$ws.Range("B24").Value = 56543
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 103021
$ws.Range("F24").Value = "Talltita"
$ws.Range("G24").Value = "Poecile montanus"
$ws.Range("H24").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("J24").ClearContents()
$ws.Range("K24").ClearContents()
$ws.Range("L24").ClearContents()
$ws.Range("M24").Value = "spel/sång"
$ws.Range("N24").ClearContents()
$ws.Range("AF24").ClearContents()
$ws.Range("Q24").Value = 554765
$ws.Range("R24").Value = 6698666
$ws.Range("Z24").ClearContents()
$ws.Range("AB24").ClearContents()

# Row 25 <= source row 24 (Id 112044195)
$ws.Range("A25").Value = 112044195
$ws.Range("B25").Value = 96348
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("J25").ClearContents()
$ws.Range("K25").ClearContents()
$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("AF25").ClearContents()
$ws.Range("Q25").Value = 554806
$ws.Range("R25").Value = 6698598
$ws.Range("Z25").ClearContents()
$ws.Range("AB25").ClearContents()

# Row 26 <= source row 25 (Id 112044171)
$ws.Range("A26").Value = 112044171
$ws.Range("B26").Value = 89686
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 658
$ws.Range("F26").Value = "Rosenticka"
$ws.Range("G26").Value = "Rhodofomes roseus"
$ws.Range("H26").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("J26").ClearContents()
$ws.Range("K26").ClearContents()
$ws.Range("L26").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("AF26").ClearContents()
$ws.Range("Q26").Value = 554758
$ws.Range("R26").Value = 6698625
$ws.Range("Z26").ClearContents()
$ws.Range("AB26").ClearContents()

# Row 27 <= source row 26 (Id 112044158)
$ws.Range("A27").Value = 112044158
$ws.Range("B27").Value = 89405
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 1202
$ws.Range("F27").Value = "Ullticka"
$ws.Range("G27").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H27").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("J27").ClearContents()
$ws.Range("K27").ClearContents()
$ws.Range("L27").ClearContents()
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("AF27").ClearContents()
$ws.Range("Q27").Value = 554756
$ws.Range("R27").Value = 6698631
$ws.Range("Z27").ClearContents()
$ws.Range("AB27").ClearContents()

# Row 28 <= source row 27 (Id 112044184)
$ws.Range("A28").Value = 112044184
$ws.Range("B28").Value = 96348
$ws.Range("D28").Value = "VU"
$ws.Range("E28").Value = 220787
$ws.Range("F28").Value = "Knärot"
$ws.Range("G28").Value = "Goodyera repens"
$ws.Range("H28").Value = "(L.) R. Br."
$ws.Range("J28").ClearContents()
$ws.Range("K28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("AF28").ClearContents()
$ws.Range("Q28").Value = 554833
$ws.Range("R28").Value = 6698646
$ws.Range("Z28").ClearContents()
$ws.Range("AB28").ClearContents()
